$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCRAdjustmentSheet_December2021")

$ws.Range("E10").Value = 47458
$ws.Range("E11").Value = 46337
$ws.Range("F11").Value = 48520
$ws.Range("E12").Value = 44587
$ws.Range("F12").Value = 46276
$ws.Range("E13").Value = 48629
$ws.Range("F13").Value = 46727
$ws.Range("E14").Value = 52946
$ws.Range("F14").Value = 52618
$ws.Range("E15").Value = 57329
$ws.Range("F15").Value = 67091
$ws.Range("E16").Value = 58479
$ws.Range("F16").Value = 67836
$ws.Range("E17").Value = 58051
$ws.Range("F17").Value = 62254
$ws.Range("E18").Value = 56464
$ws.Range("F18").Value = 59417
$ws.Range("E19").Value = 57769
$ws.Range("F19").Value = 56432
$ws.Range("E20").Value = 55272
$ws.Range("F20").Value = 57521
$ws.Range("E21").Value = 49347
$ws.Range("F21").Value = 54717
$ws.Range("F22").Value = 47717
$ws.Range("E23").Value = 28449
$ws.Range("F23").Value = 39461
$ws.Range("E24").Value = 20868
$ws.Range("F24").Value = 26354
$ws.Range("E25").Value = 16597
$ws.Range("F25").Value = 18855
$ws.Range("E26").Value = 14166
$ws.Range("F26").Value = 14284
$ws.Range("E27").Value = 16643
$ws.Range("F27").Value = 18859
$ws.Range("E28").Value = 50198
$ws.Range("F28").Value = 50230
$ws.Range("E29").Value = 48154
$ws.Range("F29").Value = 50663
$ws.Range("E30").Value = 46396
$ws.Range("F30").Value = 48080
$ws.Range("F31").Value = 47469
$ws.Range("E32").Value = 54592
$ws.Range("F32").Value = 53252
$ws.Range("E33").Value = 56246
$ws.Range("F33").Value = 67079
$ws.Range("E34").Value = 56339
$ws.Range("F34").Value = 66453
$ws.Range("E35").Value = 56854
$ws.Range("F35").Value = 61889
$ws.Range("E36").Value = 56139
$ws.Range("F36").Value = 58248
$ws.Range("E37").Value = 56108
$ws.Range("F37").Value = 56121
$ws.Range("E38").Value = 53598
$ws.Range("F38").Value = 55575
$ws.Range("E39").Value = 45954
$ws.Range("F39").Value = 52093
$ws.Range("E40").Value = 37859
$ws.Range("F40").Value = 43418
$ws.Range("E41").Value = 24901
$ws.Range("F41").Value = 34273
$ws.Range("E42").Value = 17200
$ws.Range("F42").Value = 22010
$ws.Range("E43").Value = 12651
$ws.Range("F43").Value = 14820
$ws.Range("E44").Value = 9173
$ws.Range("F44").Value = 10016
$ws.Range("E45").Value = 8333
$ws.Range("F45").Value = 9936
